# Final Completed Version of Expense Tracker v0.1
# - Re-label the old "Freelance" row as "Business Income" with an updated
#   amount/date, and add three more income sources (Salary, Youtube
#   Revenue, Interest from Savings Account) below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Income")

# Row 2: existing entry changes (Freelance -> Business Income, amount, date)
$ws.Range("A2").Value = "Business Income"
$ws.Range("B2").Value = 2000
$ws.Range("C2").Value = 45958.22928240741

# Row 3: Salary
$ws.Range("A3").Value = "Salary"
$ws.Range("B3").Value = 10000
$ws.Range("C3").Value = 45952.22928240741

# Row 4: Youtube Revenue
$ws.Range("A4").Value = "Youtube Revenue"
$ws.Range("B4").Value = 1000
$ws.Range("C4").Value = 45945.22928240741

# Row 5: Interest from Savings Account
$ws.Range("A5").Value = "Interest from Savings Account"
$ws.Range("B5").Value = 1200
$ws.Range("C5").Value = 45813.22928240741

# The new date cells (C3:C5) need the same short-date number format that
# C2 already carries (style index 1 / numFmtId 14 in the original file).
# Copying C2's formatting onto them reuses that existing style instead of
# minting a brand-new one.
$xlPasteFormats = -4122
$ws.Range("C2").Copy()
$ws.Range("C3:C5").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
